$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1665.2632
$ws.Range("I28").Value = 1647.6666
$ws.Range("J28").Value = 1731.25
$ws.Range("K28").Value = 1647.6666
$ws.Range("L28").Value = 1731.25
$ws.Range("M28").Value = -1162.6666
$ws.Range("N28").Value = -2701.25

$ws.Range("I37").Value = 1000
$ws.Range("K37").Value = 3000
$ws.Range("M37").Value = -2874

$ws.Range("H40").Value = 1817.6774
$ws.Range("I40").Value = 1558.7059
$ws.Range("K40").Value = 1558.7059
$ws.Range("M40").Value = -1383.7059

$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H64").Value = 3853.3333
$ws.Range("I64").Value = 3633.75
$ws.Range("J64").Value = 4172.727
$ws.Range("K64").Value = 3633.75
$ws.Range("L64").Value = 4172.727
$ws.Range("M64").Value = -3385.75
$ws.Range("N64").Value = -4668.727

$ws.Range("H67").Value = 3853.3333
$ws.Range("I67").Value = 3633.75
$ws.Range("J67").Value = 4172.727
$ws.Range("K67").Value = 3633.75
$ws.Range("L67").Value = 4172.727
$ws.Range("M67").Value = -2775.75
$ws.Range("N67").Value = -5888.727

$ws.Range("H82").Value = 769.0769
$ws.Range("I82").Value = 769.0769
$ws.Range("K82").Value = 2307.2307
$ws.Range("M82").Value = -1901.2307

$ws.Range("H85").Value = 769.0769
$ws.Range("I85").Value = 769.0769
$ws.Range("K85").Value = 2307.2307
$ws.Range("M85").Value = -903.2307000000001

$ws.Range("H88").Value = 2957.5454
$ws.Range("I88").Value = 2183
$ws.Range("J88").Value = 3319
$ws.Range("K88").Value = 2183
$ws.Range("L88").Value = 3319
$ws.Range("M88").Value = -1777
$ws.Range("N88").Value = -4131

$ws.Range("H91").Value = 2957.5454
$ws.Range("I91").Value = 2183
$ws.Range("J91").Value = 3319
$ws.Range("K91").Value = 2183
$ws.Range("L91").Value = 3319
$ws.Range("M91").Value = -779
$ws.Range("N91").Value = -6127

$ws.Range("H100").Value = 2081.0667
$ws.Range("I100").Value = 1625.5
$ws.Range("J100").Value = 2601.7144
$ws.Range("K100").Value = 1625.5
$ws.Range("L100").Value = 2601.7144
$ws.Range("M100").Value = -1084.5
$ws.Range("N100").Value = -3683.7144

$ws.Range("H129").Value = 938.4474
$ws.Range("I129").Value = 535.0833
$ws.Range("J129").Value = 1124.6154
$ws.Range("K129").Value = 1605.2499
$ws.Range("L129").Value = 3373.8462
$ws.Range("M129").Value = 3394.7501
$ws.Range("N129").Value = -13373.8462

$ws.Range("H130").Value = 49960
$ws.Range("J130").Value = 49960
$ws.Range("L130").Value = 49960
$ws.Range("N130").Value = -60000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2250.875
$ws.Range("I2").Value = 2348.6667
$ws.Range("J2").Value = 2087.889
$ws.Range("K2").Value = 2348.6667
$ws.Range("L2").Value = 2087.889
$ws.Range("M2").Value = -2235.6667
$ws.Range("N2").Value = -2313.889

$ws.Range("H92").Value = 28000
$ws.Range("J92").Value = 28000
$ws.Range("L92").Value = 28000
$ws.Range("N92").Value = -32992

$ws.Range("H102").Value = 2038
$ws.Range("I102").Value = 2047.5
$ws.Range("K102").Value = 2047.5
$ws.Range("M102").Value = -425.5

$ws.Range("H116").Value = 2250.875
$ws.Range("I116").Value = 2348.6667
$ws.Range("J116").Value = 2087.889
$ws.Range("K116").Value = 2348.6667
$ws.Range("L116").Value = 2087.889
$ws.Range("M116").Value = -54.66670000000022
$ws.Range("N116").Value = -6675.889

$ws.Range("H135").Value = 28999.572
$ws.Range("J135").Value = 28999.572
$ws.Range("L135").Value = 28999.572
$ws.Range("N135").Value = -39139.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2250.875
$ws.Range("I3").Value = 2348.6667
$ws.Range("J3").Value = 2087.889
$ws.Range("K3").Value = 2348.6667
$ws.Range("L3").Value = 2087.889
$ws.Range("M3").Value = -2234.6667
$ws.Range("N3").Value = -2315.889

$ws.Range("H86").Value = 1532.8928
$ws.Range("I86").Value = 1417.2778
$ws.Range("J86").Value = 1741
$ws.Range("K86").Value = 1417.2778
$ws.Range("L86").Value = 1741
$ws.Range("M86").Value = -294.2778000000001
$ws.Range("N86").Value = -3987

$ws.Range("H89").Value = 1532.8928
$ws.Range("I89").Value = 1417.2778
$ws.Range("J89").Value = 1741
$ws.Range("K89").Value = 7086.389
$ws.Range("L89").Value = 8705
$ws.Range("M89").Value = -1470.389
$ws.Range("N89").Value = -19937

$ws.Range("H94").Value = 683.7143
$ws.Range("I94").Value = 641.7
$ws.Range("J94").Value = 788.75
$ws.Range("K94").Value = 641.7
$ws.Range("L94").Value = 788.75
$ws.Range("M94").Value = -190.7
$ws.Range("N94").Value = -1690.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 850.125
$ws.Range("I16").Value = 800.1667
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 800.1667
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -513.1667
$ws.Range("N16").Value = -1574

$ws.Range("H113").Value = 850.125
$ws.Range("I113").Value = 800.1667
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 800.1667
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1369.8333
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2351.25
$ws.Range("I129").Value = 1858.0714
$ws.Range("J129").Value = 2616.8076
$ws.Range("K129").Value = 5574.2142
$ws.Range("L129").Value = 7850.4228
$ws.Range("M129").Value = -574.2142000000003
$ws.Range("N129").Value = -17850.4228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2647.3157
$ws.Range("I80").Value = 2600
$ws.Range("J80").Value = 3049.5
$ws.Range("K80").Value = 2600
$ws.Range("L80").Value = 3049.5
$ws.Range("M80").Value = -1602
$ws.Range("N80").Value = -5045.5

$ws.Range("H83").Value = 2647.3157
$ws.Range("I83").Value = 2600
$ws.Range("J83").Value = 3049.5
$ws.Range("K83").Value = 13000
$ws.Range("L83").Value = 15247.5
$ws.Range("M83").Value = -8008
$ws.Range("N83").Value = -25231.5

$ws.Range("H113").Value = 2485.8462
$ws.Range("I113").Value = 1250
$ws.Range("K113").Value = 1250
$ws.Range("M113").Value = 920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2652.4
$ws.Range("J61").Value = 4000
$ws.Range("L61").Value = 4000
$ws.Range("N61").Value = -4404

$ws.Range("H113").Value = 2652.4
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -8340

$ws.Range("H127").Value = 33382.168
$ws.Range("J127").Value = 33382.168
$ws.Range("L127").Value = 33382.168
$ws.Range("N127").Value = -43302.168

$ws.Range("H128").Value = 33999.89
$ws.Range("J128").Value = 33999.89
$ws.Range("L128").Value = 33999.89
$ws.Range("N128").Value = -43959.89

$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2711
$ws.Range("I62").Value = 2720
$ws.Range("J62").Value = 2702
$ws.Range("K62").Value = 2720
$ws.Range("L62").Value = 2702
$ws.Range("M62").Value = -2096
$ws.Range("N62").Value = -3950

$ws.Range("H65").Value = 2711
$ws.Range("I65").Value = 2720
$ws.Range("J65").Value = 2702
$ws.Range("K65").Value = 13600
$ws.Range("L65").Value = 13510
$ws.Range("M65").Value = -10480
$ws.Range("N65").Value = -19750

$ws.Range("H81").Value = 1353.8462
$ws.Range("I81").Value = 1250
$ws.Range("J81").Value = 1587.5
$ws.Range("K81").Value = 2500
$ws.Range("L81").Value = 3175
$ws.Range("M81").Value = -1439
$ws.Range("N81").Value = -5297

$ws.Range("H84").Value = 1353.8462
$ws.Range("I84").Value = 1250
$ws.Range("J84").Value = 1587.5
$ws.Range("K84").Value = 12500
$ws.Range("L84").Value = 15875
$ws.Range("M84").Value = -7196
$ws.Range("N84").Value = -26483

$ws.Range("H113").Value = 325.26666
$ws.Range("I113").Value = 189.08333
$ws.Range("J113").Value = 870
$ws.Range("K113").Value = 567.24999
$ws.Range("L113").Value = 2610
$ws.Range("M113").Value = 1602.75001
$ws.Range("N113").Value = -6950
